# Remove the "Formater ou valider: téléphone" bullet entirely (paragraph,
# its runs and its paragraph mark), per the commit "Enlegistlement:
# Waridation de téléphon" — the list item right after the
# "<<entregistrement.jsff>>" heading and right before
# "Statut doit être implanté comme une liste fixe" is dropped.

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("Formater*phone ", $false, $false, $true, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $para = $rng.Paragraphs(1)
    # Delete the whole paragraph range, including its trailing paragraph
    # mark, so the list item disappears completely instead of leaving an
    # empty bullet behind.
    $para.Range.Delete()
}
